# "Add files via upload" -- refresh the stock-screener lists in Sheet1.
# Columns: A=index, B=Buying Opportunity, C=support Zone, D=long buildup,
#          E=Short buildup, F=FII ENTERING.
# The table grows from 16 data rows (rows 2-17) to 21 data rows (rows 2-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column A's existing style (bold, centered, thin border) down to the
# newly added rows (18-22) by copying the format from an already-styled cell.
$ws.Range("A2").Copy()
$ws.Range("A18:A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    # index, B, C, D, E, F
    @(0,  "NSE:BASML",      "NSE:AARVI",      "", "NSE:AARTIIND",  "NSE:HINDPETRO"),
    @(1,  "NSE:BEPL",       "NSE:ABCAPITAL",  "", "NSE:ACC",       ""),
    @(2,  "NSE:DICIND",     "NSE:ADL",        "", "NSE:ALKEM",     ""),
    @(3,  "NSE:DIVGIITTS",  "NSE:ALKYLAMINE", "", "NSE:AMBUJACEM", ""),
    @(4,  "NSE:DMCC",       "NSE:ATAM",       "", "NSE:BOSCHLTD",  ""),
    @(5,  "NSE:ELECON",     "NSE:COFORGE",    "", "NSE:BPCL",      ""),
    @(6,  "NSE:GOKEX",      "NSE:GUJGASLTD",  "", "NSE:CANBK",     ""),
    @(7,  "NSE:GOLDETF",    "NSE:HARDWYN",    "", "NSE:CANFINHOME",""),
    @(8,  "NSE:KMSUGAR",    "NSE:IDFCFIRSTB", "", "NSE:COALINDIA", ""),
    @(9,  "NSE:MOTILALOFS", "NSE:INDOCO",     "", "NSE:COFORGE",   ""),
    @(10, "NSE:MURUDCERA",  "NSE:ITBEES",     "", "NSE:DIVISLAB",  ""),
    @(11, "NSE:NILASPACES", "NSE:ITETF",      "", "NSE:GAIL",      ""),
    @(12, "NSE:NITIRAJ",    "NSE:KAUSHALYA",  "", "NSE:GLENMARK",  ""),
    @(13, "NSE:NUVOCO",     "NSE:OBEROIRLTY", "", "NSE:GODREJPROP",""),
    @(14, "NSE:PREMEXPLN",  "NSE:RELIGARE",   "", "NSE:GRANULES",  ""),
    @(15, "NSE:RBL",        "",               "", "NSE:ICICIPRULI",""),
    @(16, "NSE:RML",        "",               "", "NSE:LTTS",      ""),
    @(17, "NSE:ROTO",       "",               "", "NSE:MCDOWELL-N",""),
    @(18, "",               "",               "", "NSE:OFSS",      ""),
    @(19, "",               "",               "", "NSE:PETRONET",  ""),
    @(20, "",               "",               "", "NSE:PIIND",     "")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
